# Updated cryptos list (Price/Volume(1h) columns) with latest figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.444.45'
$ws.Range('E2').Value = '  -1.05%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.851.83'
$ws.Range('E3').Value = '  -0.05%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.92'
$ws.Range('E5').Value = '  -0.78%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6307'
$ws.Range('E6').Value = '  -3.79%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07609'
$ws.Range('E8').Value = '  +1.38%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2974'
$ws.Range('E9').Value = '  -0.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.51'
$ws.Range('E10').Value = '  -0.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07717'
$ws.Range('E11').Value = '  +1.00%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.914.98'
$ws.Range('E12').Value = '  +3.41%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.006'
$ws.Range('E13').Value = '  -0.86%  '

$ws.Range('E14').Value = '  +0.19%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.38'
$ws.Range('E15').Value = '  -0.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009921'
$ws.Range('E16').Value = '  +2.47%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.172.27'
$ws.Range('E17').Value = '  +3.01%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.194'
$ws.Range('E18').Value = '  +1.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.556.20'
$ws.Range('E19').Value = '  -0.72%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '233.13'
$ws.Range('E20').Value = '  -1.72%  '

$ws.Range('E21').Value = '  -0.77%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.697'
$ws.Range('E22').Value = '  -0.22%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('E24').Value = '  -0.02%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.81'
$ws.Range('E25').Value = '  -2.26%  '

$ws.Range('E26').Value = '  -2.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.471'
$ws.Range('E27').Value = '  -0.78%  '

$ws.Range('E28').Value = '  -1.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.475'
$ws.Range('E29').Value = '  -1.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05780'
$ws.Range('E30').Value = '  -4.94%  '

$ws.Range('E31').Value = '  -1.26%  '

$ws.Range('E32').Value = '  -0.19%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.017'
$ws.Range('E33').Value = '  -1.47%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.891'
$ws.Range('E34').Value = '  +1.21%  '

$ws.Range('E35').Value = '  -1.22%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7224'
$ws.Range('E36').Value = '  -0.34%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.586'
$ws.Range('E37').Value = '  -0.78%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.252.86'
$ws.Range('E38').Value = '  +4.32%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.793'
$ws.Range('E39').Value = '  -0.30%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01808'
$ws.Range('E40').Value = '  +0.80%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9095'
$ws.Range('E41').Value = '  -0.26%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.092'
$ws.Range('E42').Value = '  -3.09%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.096.49'
$ws.Range('E43').Value = '  +3.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9998'
$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '67.81'
$ws.Range('E45').Value = '  +1.52%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.74'
$ws.Range('E46').Value = '  +0.54%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.322'
$ws.Range('E47').Value = '  -0.45%  '

$ws.Range('E48').Value = '  -3.49%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.174'
$ws.Range('E49').Value = '  +0.72%  '

$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4038'
$ws.Range('E50').Value = '  -0.56%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.710'
$ws.Range('E51').Value = '  +2.05%  '
